$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "HK_G_acc_G"

$values = @(
    48.141891891891895,
    48.310810810810814,
    48.479729729729733,
    47.804054054054049,
    48.310810810810814,
    48.310810810810814,
    48.986486486486484,
    49.155405405405403,
    48.648648648648653,
    48.479729729729733,
    48.986486486486484,
    49.831081081081081,
    49.662162162162161,
    50.168918918918912,
    50.168918918918912,
    50.506756756756758,
    51.013513513513509,
    50.844594594594597,
    48.986486486486484,
    49.155405405405403,
    48.986486486486484,
    47.297297297297298,
    46.95945945945946,
    46.95945945945946,
    48.479729729729733,
    48.479729729729733,
    49.324324324324323,
    50.168918918918912,
    50,
    49.493243243243242,
    47.466216216216218,
    46.95945945945946,
    47.804054054054049,
    47.635135135135137,
    48.310810810810814,
    51.520270270270274,
    47.466216216216218,
    47.635135135135137,
    47.972972972972968,
    49.155405405405403,
    48.479729729729733,
    49.831081081081081,
    49.155405405405403,
    48.141891891891895,
    48.310810810810814,
    48.141891891891895,
    51.351351351351347,
    49.324324324324323,
    49.493243243243242
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
